$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-11 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-12 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("51×36=1836", $true, $false, $false, $false, $false, $true, 1, $false, "56×52=2912", 2) | Out-Null
$d.Content.Find.Execute("26×80=2080", $true, $false, $false, $false, $false, $true, 1, $false, "83×84=6972", 2) | Out-Null
$d.Content.Find.Execute("58×43=2494", $true, $false, $false, $false, $false, $true, 1, $false, "98×64=6272", 2) | Out-Null
$d.Content.Find.Execute("72×58=4176", $true, $false, $false, $false, $false, $true, 1, $false, "15×22=330", 2) | Out-Null
$d.Content.Find.Execute("51×38=1938", $true, $false, $false, $false, $false, $true, 1, $false, "28×11=308", 2) | Out-Null
$d.Content.Find.Execute("75×85=6375", $true, $false, $false, $false, $false, $true, 1, $false, "70×74=5180", 2) | Out-Null
$d.Content.Find.Execute("84×51=4284", $true, $false, $false, $false, $false, $true, 1, $false, "59×22=1298", 2) | Out-Null
$d.Content.Find.Execute("43×25=1075", $true, $false, $false, $false, $false, $true, 1, $false, "23×20=460", 2) | Out-Null
$d.Content.Find.Execute("78×40=3120", $true, $false, $false, $false, $false, $true, 1, $false, "63×76=4788", 2) | Out-Null
$d.Content.Find.Execute("12×90=1080", $true, $false, $false, $false, $false, $true, 1, $false, "96×19=1824", 2) | Out-Null
$d.Content.Find.Execute("70×84=5880", $true, $false, $false, $false, $false, $true, 1, $false, "65×37=2405", 2) | Out-Null
$d.Content.Find.Execute("99×33=3267", $true, $false, $false, $false, $false, $true, 1, $false, "38×13=494", 2) | Out-Null
$d.Content.Find.Execute("67×83=5561", $true, $false, $false, $false, $false, $true, 1, $false, "71×62=4402", 2) | Out-Null
$d.Content.Find.Execute("19×64=1216", $true, $false, $false, $false, $false, $true, 1, $false, "24×19=456", 2) | Out-Null
$d.Content.Find.Execute("26×47=1222", $true, $false, $false, $false, $false, $true, 1, $false, "34×72=2448", 2) | Out-Null
$d.Content.Find.Execute("68×85=5780", $true, $false, $false, $false, $false, $true, 1, $false, "43×59=2537", 2) | Out-Null
$d.Content.Find.Execute("93×60=5580", $true, $false, $false, $false, $false, $true, 1, $false, "76×76=5776", 2) | Out-Null
$d.Content.Find.Execute("96×62=5952", $true, $false, $false, $false, $false, $true, 1, $false, "40×83=3320", 2) | Out-Null
$d.Content.Find.Execute("47×14=658", $true, $false, $false, $false, $false, $true, 1, $false, "82×42=3444", 2) | Out-Null
$d.Content.Find.Execute("96×86=8256", $true, $false, $false, $false, $false, $true, 1, $false, "40×31=1240", 2) | Out-Null
$d.Content.Find.Execute("39×36=1404", $true, $false, $false, $false, $false, $true, 1, $false, "85×34=2890", 2) | Out-Null
$d.Content.Find.Execute("40×81=3240", $true, $false, $false, $false, $false, $true, 1, $false, "94×99=9306", 2) | Out-Null
$d.Content.Find.Execute("87×76=6612", $true, $false, $false, $false, $false, $true, 1, $false, "46×50=2300", 2) | Out-Null
$d.Content.Find.Execute("88×73=6424", $true, $false, $false, $false, $false, $true, 1, $false, "35×11=385", 2) | Out-Null
$d.Content.Find.Execute("85×42=3570", $true, $false, $false, $false, $false, $true, 1, $false, "94×98=9212", 2) | Out-Null
